$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> ISO date string that replaces the numeric Excel
# date serial currently stored (with the custom YYYY-MM-DD HH:MM:SS style)
# in column A. The new cells become plain text, no special style.
$dates = @{
    2 = "2022-07-01"
    3 = "2022-07-01"
    4 = "2022-07-01"
    5 = "2022-07-01"
    6 = "2022-07-01"
    7 = "2022-07-01"
    8 = "2022-07-01"
    9 = "2022-07-01"
    10 = "2022-07-01"
    11 = "2022-06-30"
    12 = "2022-06-30"
    13 = "2022-06-30"
    14 = "2022-06-30"
    15 = "2022-06-30"
    16 = "2022-06-30"
    17 = "2022-06-30"
    18 = "2022-06-30"
    19 = "2022-06-30"
    20 = "2022-06-30"
    21 = "2022-06-30"
    22 = "2022-06-30"
    23 = "2022-06-30"
    24 = "2022-06-30"
    25 = "2022-06-30"
    26 = "2022-06-29"
    27 = "2022-06-29"
    28 = "2022-06-29"
    29 = "2022-06-29"
    30 = "2022-06-29"
    31 = "2022-06-29"
    32 = "2022-06-29"
    33 = "2022-06-28"
    34 = "2022-06-28"
    35 = "2022-06-28"
    36 = "2022-06-28"
    37 = "2022-06-28"
    38 = "2022-06-28"
    39 = "2022-06-27"
    40 = "2022-06-27"
    41 = "2022-06-27"
    42 = "2022-06-27"
    43 = "2022-06-25"
    44 = "2022-06-25"
    45 = "2022-06-24"
    46 = "2022-06-24"
    47 = "2022-06-24"
    48 = "2022-06-23"
    49 = "2022-06-23"
    50 = "2022-06-23"
    51 = "2022-06-23"
    52 = "2022-06-23"
    53 = "2022-06-22"
    54 = "2022-06-22"
    55 = "2022-06-21"
    56 = "2022-06-21"
    57 = "2022-06-20"
    58 = "2022-06-20"
    59 = "2022-06-20"
    60 = "2022-06-17"
    61 = "2022-06-16"
    62 = "2022-06-16"
    63 = "2022-06-15"
    64 = "2022-06-15"
    65 = "2022-06-14"
    66 = "2022-06-14"
    67 = "2022-06-13"
    68 = "2022-06-11"
    69 = "2022-06-07"
    70 = "2022-06-03"
    71 = "2022-05-31"
}

foreach ($row in ($dates.Keys | Sort-Object)) {
    $cell = $ws.Cells.Item($row, 1)
    # Leading apostrophe forces text interpretation so Excel doesn't
    # re-parse the ISO string back into a date serial number.
    $cell.Value = "'" + $dates[$row]
    # Drop back to the default "Normal" style so no explicit number
    # format / style index is left on the cell.
    $cell.Style = "Normal"
}
